# Generate Report for Handoff
#
# The three status files 8b95212a-...md, 9d8144f4-...md and 73baee7e-...md
# (rows 6, 7, 8 on every sheet) get "rotated" up by one row:
#   row6 -> row7, row7 -> row8, row8 -> row6
# and the file that ends up on row 8 (8b95212a-...md) is now fully handed
# off, so its Status flips to "Ready for handoff" and it receives a fresh
# handoff timestamp.

$wb = $excel.ActiveWorkbook

function Set-RowAndLink {
    param(
        $ws,
        [int]$row,
        [hashtable]$values,
        [hashtable]$linkDisplays
    )

    foreach ($col in $values.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $values[$col]
    }

    if ($linkDisplays.Count -gt 0) {
        foreach ($h in $ws.Hyperlinks) {
            $addr = $h.Range.Address(0, 0)
            foreach ($col in $linkDisplays.Keys) {
                if ($addr -eq "$col$row") {
                    $h.TextToDisplay = $linkDisplays[$col]
                }
            }
        }
    }
}

# ---------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

Set-RowAndLink $ws 6 @{
    "A" = "9d8144f4-0c54-448a-90fa-c0a02fb30668.md"
    "B" = "In Translation"
    "C" = "In Translation"
    "D" = "2016-03-22 22:46:29"
} @{ "A" = "9d8144f4-0c54-448a-90fa-c0a02fb30668.md" }

Set-RowAndLink $ws 7 @{
    "A" = "73baee7e-d827-4a80-9bc6-a83dc0ac0602.md"
    "B" = "Ready for handoff"
    "C" = "Ready for handoff"
    "D" = "2016-03-22 22:43:40"
} @{ "A" = "73baee7e-d827-4a80-9bc6-a83dc0ac0602.md" }

Set-RowAndLink $ws 8 @{
    "A" = "8b95212a-cc2f-4185-9ddb-738e68c91732.md"
    "B" = "Ready for handoff"
    "C" = "Ready for handoff"
    "D" = "2016-03-22 22:50:18"
} @{ "A" = "8b95212a-cc2f-4185-9ddb-738e68c91732.md" }

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

Set-RowAndLink $ws 6 @{
    "A" = "9d8144f4-0c54-448a-90fa-c0a02fb30668.md"
    "C" = "In Translation"
    "D" = "9d8144f4-0c54-448a-90fa-c0a02fb30668.8fdfc45c12a1b4ecd0c9602c5414a318d0be7c32.zh-cn.xlf"
    "E" = "2016-03-22 22:46:24"
} @{
    "A" = "9d8144f4-0c54-448a-90fa-c0a02fb30668.md"
    "D" = "9d8144f4-0c54-448a-90fa-c0a02fb30668.8fdfc45c12a1b4ecd0c9602c5414a318d0be7c32.zh-cn.xlf"
}

Set-RowAndLink $ws 7 @{
    "A" = "73baee7e-d827-4a80-9bc6-a83dc0ac0602.md"
    "C" = "Ready for handoff"
    "D" = "73baee7e-d827-4a80-9bc6-a83dc0ac0602.9d603e7a0a8b691f7b0c595a2378f538092f8fd6.zh-cn.xlf"
    "E" = "2016-03-22 22:43:37"
} @{
    "A" = "73baee7e-d827-4a80-9bc6-a83dc0ac0602.md"
    "D" = "73baee7e-d827-4a80-9bc6-a83dc0ac0602.9d603e7a0a8b691f7b0c595a2378f538092f8fd6.zh-cn.xlf"
}

Set-RowAndLink $ws 8 @{
    "A" = "8b95212a-cc2f-4185-9ddb-738e68c91732.md"
    "C" = "Ready for handoff"
    "D" = "8b95212a-cc2f-4185-9ddb-738e68c91732.2219b3df99a62afc43ba998e239bf69c985229f5.zh-cn.xlf"
    "E" = "2016-03-22 22:50:14"
} @{
    "A" = "8b95212a-cc2f-4185-9ddb-738e68c91732.md"
    "D" = "8b95212a-cc2f-4185-9ddb-738e68c91732.2219b3df99a62afc43ba998e239bf69c985229f5.zh-cn.xlf"
}

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

Set-RowAndLink $ws 6 @{
    "A" = "9d8144f4-0c54-448a-90fa-c0a02fb30668.md"
    "C" = "In Translation"
    "D" = "9d8144f4-0c54-448a-90fa-c0a02fb30668.8fdfc45c12a1b4ecd0c9602c5414a318d0be7c32.de-de.xlf"
    "E" = "2016-03-22 22:46:29"
} @{
    "A" = "9d8144f4-0c54-448a-90fa-c0a02fb30668.md"
    "D" = "9d8144f4-0c54-448a-90fa-c0a02fb30668.8fdfc45c12a1b4ecd0c9602c5414a318d0be7c32.de-de.xlf"
}

Set-RowAndLink $ws 7 @{
    "A" = "73baee7e-d827-4a80-9bc6-a83dc0ac0602.md"
    "C" = "Ready for handoff"
    "D" = "73baee7e-d827-4a80-9bc6-a83dc0ac0602.9d603e7a0a8b691f7b0c595a2378f538092f8fd6.de-de.xlf"
    "E" = "2016-03-22 22:43:40"
} @{
    "A" = "73baee7e-d827-4a80-9bc6-a83dc0ac0602.md"
    "D" = "73baee7e-d827-4a80-9bc6-a83dc0ac0602.9d603e7a0a8b691f7b0c595a2378f538092f8fd6.de-de.xlf"
}

Set-RowAndLink $ws 8 @{
    "A" = "8b95212a-cc2f-4185-9ddb-738e68c91732.md"
    "C" = "Ready for handoff"
    "D" = "8b95212a-cc2f-4185-9ddb-738e68c91732.2219b3df99a62afc43ba998e239bf69c985229f5.de-de.xlf"
    "E" = "2016-03-22 22:50:18"
} @{
    "A" = "8b95212a-cc2f-4185-9ddb-738e68c91732.md"
    "D" = "8b95212a-cc2f-4185-9ddb-738e68c91732.2219b3df99a62afc43ba998e239bf69c985229f5.de-de.xlf"
}
